$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting the existing rows 42:62 down to 43:63
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with the new data record
$ws.Range("A42").Value = 4
$ws.Range("B42").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C42").Value = "Los Lagos"
$ws.Range("D42").NumberFormat = $ws.Range("D43").NumberFormat
$ws.Range("D42").Value = 44574
$ws.Range("E42").Value = 10
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100103
$ws.Range("H42").Value = "Frutos de hueso (carozo)"
$ws.Range("I42").Value = 100103001
$ws.Range("J42").Value = "Cereza"
$ws.Range("K42").Value = "Santina"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 7500
$ws.Range("O42").Value = 8000
$ws.Range("P42").Value = 7750
$ws.Range("Q42").Value = "$/bandeja 10 kilos"
$ws.Range("R42").Value = "Provincia de Curicó"
$ws.Range("S42").Value = 775
$ws.Range("T42").Value = 10
